# Apply cell content updates to Sheet1 as described by the commit diff.
# The business_class / class labels are being replaced with more
# descriptive human-readable (Chinese) short names ("extract by zuh eid").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "business_short"
$ws.Range("C2").Value = "登录按钮"
$ws.Range("C3").Value = "旅客检索页面加载"
$ws.Range("C4").Value = "用户注销按钮"
$ws.Range("C5").Value = "登录按钮"
$ws.Range("C6").Value = "旅客检索页面加载"
$ws.Range("C7").Value = "用户临时退出页面加载"
$ws.Range("C8").Value = "登录按钮"
$ws.Range("C9").Value = "旅客检索页面加载"
$ws.Range("C10").Value = "用户注销按钮"
$ws.Range("C11").Value = "登录按钮"
$ws.Range("C12").Value = "旅客检索页面加载"
$ws.Range("C13").Value = "旅客检索清空按钮"
$ws.Range("C14").Value = "旅客检索清空按钮"
$ws.Range("C15").Value = "登录按钮"
$ws.Range("C16").Value = "旅客检索页面加载"
$ws.Range("C17").Value = "用户注销按钮"
$ws.Range("C18").Value = "登录按钮"
$ws.Range("C19").Value = "旅客检索页面加载"
$ws.Range("C20").Value = "用户注销按钮"
$ws.Range("C21").Value = "登录按钮"
$ws.Range("C22").Value = "旅客检索页面加载"
$ws.Range("C23").Value = "用户注销按钮"
$ws.Range("C24").Value = "登录按钮"
$ws.Range("C25").Value = "旅客检索页面加载"
$ws.Range("C26").Value = "用户注销按钮"
$ws.Range("C27").Value = "登录按钮"
$ws.Range("C33").Value = "用户注销按钮"
$ws.Range("C34").Value = "登录按钮"
$ws.Range("C35").Value = "旅客检索页面加载"
$ws.Range("C36").Value = "用户注销按钮"
$ws.Range("C37").Value = "登录按钮"
$ws.Range("C38").Value = "旅客检索页面加载"
$ws.Range("C39").Value = "用户注销按钮"
$ws.Range("C40").Value = "用户注销按钮"
$ws.Range("C41").Value = "登录按钮"
$ws.Range("C42").Value = "旅客检索页面加载"
$ws.Range("C43").Value = "用户注销按钮"
$ws.Range("C44").Value = "用户注销按钮"
$ws.Range("C45").Value = "登录按钮"
$ws.Range("C46").Value = "旅客检索页面加载"
$ws.Range("C47").Value = "用户注销按钮"
$ws.Range("C48").Value = "登录按钮"
$ws.Range("C49").Value = "旅客检索页面加载"
$ws.Range("C50").Value = "登录按钮"
$ws.Range("C51").Value = "旅客检索页面加载"
$ws.Range("C52").Value = "用户注销按钮"
$ws.Range("C53").Value = "登录按钮"
$ws.Range("C54").Value = "登录按钮"
$ws.Range("C56").Value = "旅客检索页面加载"
$ws.Range("C57").Value = "登录按钮"
$ws.Range("C58").Value = "旅客检索页面加载"
$ws.Range("C59").Value = "用户注销按钮"
$ws.Range("C60").Value = "登录按钮"
$ws.Range("C62").Value = "用户注销按钮"
$ws.Range("C63").Value = "登录按钮"
$ws.Range("C64").Value = "旅客检索页面加载"
$ws.Range("C65").Value = "用户注销按钮"
$ws.Range("C66").Value = "登录按钮"
$ws.Range("C70").Value = "旅客详情刷新按钮"
$ws.Range("E70").Value = "旅客详情"
$ws.Range("C71").Value = "旅客详情刷新按钮"
$ws.Range("E71").Value = "旅客详情"
$ws.Range("C75").Value = "列表切换旅客"
$ws.Range("C76").Value = "旅客详情刷新按钮"
$ws.Range("E76").Value = "旅客详情"
$ws.Range("C77").Value = "旅客详情刷新按钮"
$ws.Range("E77").Value = "旅客详情"
$ws.Range("C78").Value = "旅客详情刷新按钮"
$ws.Range("E78").Value = "旅客详情"
$ws.Range("C79").Value = "旅客详情刷新按钮"
$ws.Range("E79").Value = "旅客详情"
$ws.Range("C80").Value = "旅客详情刷新按钮"
$ws.Range("E80").Value = "旅客详情"
$ws.Range("C81").Value = "旅客详情刷新按钮"
$ws.Range("E81").Value = "旅客详情"
$ws.Range("C82").Value = "旅客详情刷新按钮"
$ws.Range("E82").Value = "旅客详情"
$ws.Range("C83").Value = "旅客详情刷新按钮"
$ws.Range("E83").Value = "旅客详情"
$ws.Range("C84").Value = "列表切换旅客"
$ws.Range("C85").Value = "旅客详情刷新按钮"
$ws.Range("E85").Value = "旅客详情"
$ws.Range("C86").Value = "旅客详情刷新按钮"
$ws.Range("E86").Value = "旅客详情"
$ws.Range("C87").Value = "旅客详情刷新按钮"
$ws.Range("E87").Value = "旅客详情"
$ws.Range("C89").Value = "旅客详情刷新按钮"
$ws.Range("E89").Value = "旅客详情"
$ws.Range("C90").Value = "旅客详情刷新按钮"
$ws.Range("E90").Value = "旅客详情"
$ws.Range("C91").Value = "旅客详情刷新按钮"
$ws.Range("E91").Value = "旅客详情"
$ws.Range("C92").Value = "旅客详情刷新按钮"
$ws.Range("E92").Value = "旅客详情"
$ws.Range("C93").Value = "旅客检索页面加载"
$ws.Range("C94").Value = "用户注销按钮"
$ws.Range("C95").Value = "登录按钮"
$ws.Range("C97").Value = "用户注销按钮"
$ws.Range("C98").Value = "登录按钮"
$ws.Range("C99").Value = "旅客检索页面加载"
$ws.Range("C100").Value = "用户注销按钮"
$ws.Range("C101").Value = "登录按钮"
$ws.Range("C102").Value = "旅客检索页面加载"
$ws.Range("C103").Value = "已登机查询按钮"
$ws.Range("E103").Value = "已/未登机查询入口"
$ws.Range("C104").Value = "登录按钮"
$ws.Range("C105").Value = "旅客检索页面加载"
$ws.Range("C106").Value = "旅客提取页签检索按钮"
$ws.Range("C107").Value = "重打登机牌按钮"
$ws.Range("E107").Value = "重打牌登机牌"
$ws.Range("C108").Value = "重打登机牌按钮"
$ws.Range("E108").Value = "重打牌登机牌"
$ws.Range("C109").Value = "旅客提取页签检索按钮"
$ws.Range("C110").Value = "序号输入框选中旅客"
$ws.Range("C111").Value = "重打登机牌按钮"
$ws.Range("E111").Value = "重打牌登机牌"
$ws.Range("C112").Value = "重打登机牌按钮"
$ws.Range("E112").Value = "重打牌登机牌"
$ws.Range("C113").Value = "旅客提取页签检索按钮"
$ws.Range("C114").Value = "序号输入框选中旅客"
$ws.Range("C115").Value = "重打登机牌按钮"
$ws.Range("E115").Value = "重打牌登机牌"
$ws.Range("C116").Value = "重打登机牌按钮"
$ws.Range("E116").Value = "重打牌登机牌"
